$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.929.98"
$ws.Range("E2").Value = "  -3.93%  "
$ws.Range("D3").Value = "2.450.75"
$ws.Range("E3").Value = "  -3.10%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'309.46"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'92.49"
$ws.Range("E6").Value = "  -7.54%  "
$ws.Range("D7").Value = "'0.549"
$ws.Range("E7").Value = "  -3.10%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.492"
$ws.Range("E9").Value = "  -5.51%  "
$ws.Range("D10").Value = "'33.01"
$ws.Range("E10").Value = "  -7.25%  "
$ws.Range("D11").Value = "'0.0771"
$ws.Range("E11").Value = "  -3.83%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.92"
$ws.Range("E13").Value = "  -5.77%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.829.79"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").Value = "2.445.52"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "'14.69"
$ws.Range("E16").Value = "  -3.73%  "
$ws.Range("D17").Value = "'0.773"
$ws.Range("E17").Value = "  -4.35%  "
$ws.Range("D18").Value = "40.897.79"
$ws.Range("E18").Value = "  -3.97%  "
$ws.Range("D19").Value = "'6.23"
$ws.Range("E19").Value = "  -6.94%  "
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("D21").Value = "'10.97"
$ws.Range("E21").Value = "  -9.81%  "
$ws.Range("D22").Value = "'67.60"
$ws.Range("E22").Value = "  -2.63%  "
$ws.Range("D23").Value = "'233.63"
$ws.Range("E23").Value = "  -3.71%  "
$ws.Range("D24").Value = "'2.73"
$ws.Range("E24").Value = "  -4.47%  "
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("D26").Value = "'1.88"
$ws.Range("E26").Value = "  -7.51%  "
$ws.Range("D27").Value = "'23.64"
$ws.Range("E27").Value = "  -7.20%  "
$ws.Range("E28").Value = "  -5.98%  "
$ws.Range("D29").Value = "'9.50"
$ws.Range("E29").Value = "  -6.02%  "
$ws.Range("D30").Value = "'35.49"
$ws.Range("E30").Value = "  -7.12%  "
$ws.Range("D31").Value = "'150.59"
$ws.Range("E31").Value = "  -4.34%  "
$ws.Range("D32").Value = "'5.42"
$ws.Range("E32").Value = "  -5.22%  "
$ws.Range("E33").Value = "  -5.73%  "
$ws.Range("E34").Value = "  -3.90%  "
$ws.Range("D35").Value = "'0.0733"
$ws.Range("E35").Value = "  -6.14%  "
$ws.Range("D36").Value = "'2.96"
$ws.Range("E36").Value = "  -5.47%  "
$ws.Range("D37").Value = "'16.65"
$ws.Range("E37").Value = "  -7.21%  "
$ws.Range("D38").Value = "'1.84"
$ws.Range("E38").Value = "  -6.42%  "
$ws.Range("E39").Value = "  -3.74%  "
$ws.Range("D40").Value = "'0.101"
$ws.Range("E40").Value = "  -8.51%  "
$ws.Range("D41").Value = "'4.13"
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'19.71"
$ws.Range("E43").Value = "  -11.97%  "
$ws.Range("D44").Value = "1.956.32"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").Value = "'0.0281"
$ws.Range("E45").Value = "  -6.03%  "
$ws.Range("D46").Value = "'2.99"
$ws.Range("E46").Value = "  -8.27%  "
$ws.Range("D47").Value = "'8.52"
$ws.Range("E47").Value = "  -4.13%  "
$ws.Range("D48").Value = "'68.71"
$ws.Range("E48").Value = "  -4.62%  "
$ws.Range("D49").Value = "'95.88"
$ws.Range("E49").Value = "  -4.76%  "
$ws.Range("D50").Value = "'0.175"
$ws.Range("E50").Value = "  -6.80%  "
$ws.Range("D51").Value = "'73.51"
$ws.Range("E51").Value = "  -7.09%  "

# Reset number format footprint on cells that needed text-coercion
# (values that would otherwise be auto-parsed as numbers)
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
